$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-02-20 Thursday" "2025-02-21 Friday"

Replace-Text "675×8=5400" "970×4=3880"
Replace-Text "702×5=3510" "493×7=3451"
Replace-Text "947×7=6629" "171×2=342"
Replace-Text "195×2=390" "367×3=1101"
Replace-Text "831×3=2493" "572×7=4004"

Replace-Text "496×6=2976" "987×4=3948"
Replace-Text "438×5=2190" "981×3=2943"
Replace-Text "523×8=4184" "117×9=1053"
Replace-Text "935×8=7480" "588×5=2940"
Replace-Text "510×8=4080" "218×5=1090"

Replace-Text "468×2=936" "471×7=3297"
Replace-Text "356×6=2136" "491×8=3928"
Replace-Text "226×7=1582" "548×3=1644"
Replace-Text "420×5=2100" "210×7=1470"
Replace-Text "822×6=4932" "938×5=4690"

Replace-Text "316×9=2844" "527×3=1581"
Replace-Text "462×7=3234" "693×4=2772"
Replace-Text "944×8=7552" "276×8=2208"
Replace-Text "391×3=1173" "187×2=374"
Replace-Text "891×4=3564" "336×9=3024"

Replace-Text "754×3=2262" "610×6=3660"
Replace-Text "528×2=1056" "523×4=2092"
Replace-Text "301×5=1505" "379×3=1137"
Replace-Text "235×3=705" "712×9=6408"
Replace-Text "259×4=1036" "596×8=4768"
